$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39; existing rows 39-60 shift down to 40-61.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44784
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112013
$ws.Range("G39").Value = "Alcachofa"
$ws.Range("H39").Value = "Española"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 14000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 14400
$ws.Range("N39").Value = "`$/caja 30 unidades"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 480
$ws.Range("Q39").Value = 30
$ws.Range("R39").Value = "Hortaliza"
